# Journal de travail - Diogo
# Update the work journal: row 48-49 time adjustments, add a new work entry
# (row 52-53) describing the creation of the automatic installation script.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 48: time worked on "documentation d'installation" entry increased from 1.5 to 2
$ws.Range("D48").Value = 2

# Row 49: time worked on "vérifier l'ergonomie" entry increased from 1 to 2
$ws.Range("D49").Value = 2

# Row 52: new entry - "documentation d'installation" task continued, 2 hours
$ws.Range("B52").Value = "documentation d'installation"
$ws.Range("D52").Value = 2

# Row 53: new entry - creation of the automatic installation script, 2 hours
$ws.Range("B53").Value = "création du script d'installation automatique"
$ws.Range("D53").Value = 2

# Reflect where the user was last working when they saved the file
$ws.Range("E53").Select()
